# Change filtering common results by website title instead of urls.
# - DuckDuckGo sheet's A3 value ("privacy") is replaced with "quality assurance"
#   (the same value used on the Google sheet), since common results are now
#   matched by website title instead of URL.
# - Update the current selection on each sheet (no longer referencing column C/D).

$wb = $excel.ActiveWorkbook

$wsGoogle = $wb.Worksheets.Item("Google")
$wsDuckDuckGo = $wb.Worksheets.Item("DuckDuckGo")

# Update the value that used to be "privacy" to "quality assurance"
$wsDuckDuckGo.Range("A3").Value = "quality assurance"

# Update selections on each sheet
$wsGoogle.Range("A7").Select()
$wsDuckDuckGo.Range("A9").Select()
